$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation (2026/01/24, Sat, 16:00, rank 201) was recorded for a date
# that had already been appended at the bottom of the sheet (2026/12/29 ...).
# It belongs right before the existing row 704, so insert a blank row there
# and push everything from 704:745 down to 705:746 (dimension grows to D746).
$ws.Rows.Item(704).Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), not real Excel date
# serials. A direct .Value assignment of a date-shaped string gets
# auto-converted to a date number, so force it in as text via a leading
# apostrophe, then reset the style back to the sheet's default (no quote
# prefix / no number format) to match the rest of the column.
$ws.Cells.Item(704, 1).Value = "'2026/01/24"
$ws.Cells.Item(704, 1).Style = "Normal"

$ws.Cells.Item(704, 2).Value = "土"
$ws.Cells.Item(704, 3).Value = 16
$ws.Cells.Item(704, 4).Value = 201
